$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.401.40'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.849.03'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.33'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6299'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07625'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2938'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.48'
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07748'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '1.848.43'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.010'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001084'
$ws.Range('E14').Value = '  +8.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6793'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.52'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '2.088.21'
$ws.Range('E17').Value = '  -7.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.146'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '29.431.49'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.63'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.454'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.16'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.375'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.64'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.306'
$ws.Range('E30').Value = '  +4.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05631'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.115'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.047'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.780'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.231.53'
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('E40').Value = '  -1.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.480'
$ws.Range('E41').Value = '  +4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9070'
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.37'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.09'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000121'
$ws.Range('E46').Value = '  +4.31%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.166'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4012'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.018'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1123'
$ws.Range('E51').Value = '  -0.70%  '
